# Updated slides for second training
# - Bump the cached "today" date stamp on the master (date field text)
# - Bump the cover-slide date
# - Bump the three week-day dates on the "Zeitplan" slide
# - Bump the GitLab repo URL used on the "Organisation" slide and drop the
#   now-redundant blank paragraph that used to separate it from the final
#   bullet paragraph

function Trim-Para([string]$s) {
    return $s.TrimEnd("`r", "`n")
}

$p = $ppt.ActivePresentation

# 1) Slide master: literal cached text of the date field ("Rectangle 6")
#    19.06.2024 -> 06.07.2024
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ((Trim-Para $shp.TextFrame.TextRange.Text) -eq "19.06.2024") {
            $shp.TextFrame.TextRange.Text = "06.07.2024"
        }
    }
}

# 2) Slide 1 (title slide): "17.06.2024, Daniel Krämer" -> "08.07.2024, ..."
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ((Trim-Para $tr.Text) -eq "17.06.2024, Daniel Krämer") {
            $tr.Paragraphs(1, 1).Runs(1, 1).Text = "08.07.2024, Daniel Krämer"
        }
    }
}

# 3) Slide 9 ("Zeitplan"): the three weekday/date lines
$s9 = $p.Slides.Item(9)
for ($i = 1; $i -le $s9.Shapes.Count; $i++) {
    $shp = $s9.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j, 1)
            $txt = Trim-Para $para.Text
            switch ($txt) {
                "Montag, 17.06.2024"   { $para.Runs(1, 1).Text = "Montag, 08.07.2024" }
                "Dienstag, 18.06.2024" { $para.Runs(1, 1).Text = "Dienstag, 09.07.2024" }
                "Mittwoch, 19.06.2024" { $para.Runs(1, 1).Text = "Mittwoch, 10.07.2024" }
            }
        }
    }
}

# 4) Slide 10 ("Organisation"): bump the GitLab repo URL and remove the
#    blank paragraph that used to follow it.
$s10 = $p.Slides.Item(10)
for ($i = 1; $i -le $s10.Shapes.Count; $i++) {
    $shp = $s10.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*gitlab-24.06*") {
            $paraCount = $tr.Paragraphs().Count
            for ($j = $paraCount; $j -ge 1; $j--) {
                $para = $tr.Paragraphs($j, 1)
                $txt = Trim-Para $para.Text
                if ($txt -eq "https://github.com/anderscore-gmbh/gitlab-24.06") {
                    $para.Runs(1, 1).Text = "https://github.com/anderscore-gmbh/gitlab-24.07"
                    # The empty paragraph right after the link (marL=0, lvl=1,
                    # buNone, no bullet) is no longer needed; drop it.
                    $next = $tr.Paragraphs($j + 1, 1)
                    $nextTxt = Trim-Para $next.Text
                    if ($nextTxt -eq "" -and $next.ParagraphFormat.Bullet.Visible -eq 0) {
                        $next.Delete()
                    }
                }
            }
        }
    }
}
